$p = $ppt.ActivePresentation

# --- Add a new blank slide as slide 2 (after the existing title slide) ---
[void]$p.Slides.Add(2, 12)

# --- Merge the two runs "Cliente " + "(Navegador Web)" into a single run ---
$s1 = $p.Slides.Item(1)
$grp = $s1.Shapes.Item(3)
$clienteShape = $grp.GroupItems.Item(2)
# Force a real text replace (same value is a no-op) so the run split collapses.
$clienteShape.TextFrame.TextRange.Text = "zzz"
$clienteShape.TextFrame.TextRange.Text = "Cliente (Navegador Web)"
